$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-09-12"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 09-12)"

# Update September total (row 10) and yearly Total (row 14) for the "2022" column (I)
$ws.Range("I10").Value = 61
$ws.Range("I14").Value = 1198
